$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "Attributes.Boolean"
$ws.Range("F2").Value = $true

# Widen the new column to roughly match the author's manual resize.
$ws.Columns.Item(6).ColumnWidth = 19.8

[void]$ws.Range("C9").Select()
